$wb = $excel.ActiveWorkbook

# Step 1: clear all text cells in column A (and header B1) across all three sheets so the
# existing shared-string table entries become unreferenced and get garbage collected on save.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1:B7").ClearContents()
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1:B26").ClearContents()
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A1:B30").ClearContents()

# Step 2: re-write every cell value in the exact order the strings should appear in the
# rebuilt shared-string table (matches the order used by the target workbook).
$ws1.Range("A1").Value = "Item Name"
$ws2.Range("A1").Value = "Item Name"
$ws3.Range("A1").Value = "Item Name"
$ws1.Range("B1").Value = "Quantity"
$ws2.Range("B1").Value = "Quantity"
$ws3.Range("B1").Value = "Quantity"
$ws2.Range("A2").Value = "Borewell Developing Charges"
$ws3.Range("A2").Value = "Borewell Developing Charges"
$ws2.Range("A4").Value = "Supply Fitting of ISI UPVC Pipe 40mm dia"
$ws2.Range("A5").Value = "Supply Fitting of Nylon Rope 14mm thick"
$ws3.Range("A10").Value = "Supply Fitting of Nylon Rope 14mm thick"
$ws2.Range("A7").Value = "Pump erection charge upto 3hp"
$ws2.Range("A8").Value = "Supply fitting of NRV 40mm"
$ws3.Range("A13").Value = "Supply fitting of NRV 40mm"
$ws2.Range("A9").Value = "Supply Fitting of SS Adaptor 40mm dia"
$ws3.Range("A14").Value = "Supply Fitting of SS Adaptor 40mm dia"
$ws2.Range("A10").Value = "Fabrication and Installation of Well Protection Cover"
$ws3.Range("A12").Value = "Fabrication and Installation of Well Protection Cover"
$ws2.Range("A11").Value = "Fabrication and Installation of Metallic pump house 1.2 x 0.9 x 2m"
$ws3.Range("A16").Value = "Fabrication and Installation of Metallic pump house 1.2 x 0.9 x 2m"
$ws2.Range("A12").Value = "Fabrication and Installation of Panel board Single phase"
$ws3.Range("A17").Value = "Fabrication and Installation of Panel board Single phase"
$ws2.Range("A13").Value = "Fabrication and Installation of Iron Structure 4m for 3000 ltr water tank"
$ws3.Range("A18").Value = "Fabrication and Installation of Iron Structure 4m for 3000 ltr water tank"
$ws2.Range("A14").Value = "Supply fitting of ISI PVC Water Tank"
$ws3.Range("A19").Value = "Supply fitting of ISI PVC Water Tank"
$ws2.Range("A15").Value = "Tank connector 50mm"
$ws3.Range("A20").Value = "Tank connector 50mm"
$ws2.Range("A16").Value = "Tank connector 40mm"
$ws2.Range("A17").Value = "Supply fitting of PVC ball valve 50mm"
$ws3.Range("A21").Value = "Supply fitting of PVC ball valve 50mm"
$ws2.Range("A18").Value = "Supply fitting of PVC ball valve 40mm"
$ws3.Range("A22").Value = "Supply fitting of PVC ball valve 40mm"
$ws2.Range("A19").Value = "Providing and fixing ISI PVC Pipes includes jointing of pipes with one step PVC solvent cement, trenching, refilling & testing of joints complete as per direction of Engineer in Charge 50mm 10kg/cm2."
$ws3.Range("A23").Value = "Providing and fixing ISI PVC Pipes includes jointing of pipes with one step PVC solvent cement, trenching, refilling & testing of joints complete as per direction of Engineer in Charge 50mm 10kg/cm2."
$ws2.Range("A20").Value = "Providing and fixing ISI PVC Pipes includes jointing of pipes with one step PVC solvent cement, trenching, refilling & testing of joints complete as per direction of Engineer in Charge 40mm 10kg/cm2."
$ws3.Range("A24").Value = "Providing and fixing ISI PVC Pipes includes jointing of pipes with one step PVC solvent cement, trenching, refilling & testing of joints complete as per direction of Engineer in Charge 40mm 10kg/cm2."
$ws2.Range("A21").Value = "Providing and fixing ISI PVC Pipes includes jointing of pipes with one step PVC solvent cement, trenching, refilling & testing of joints complete as per direction of Engineer in Charge 32mm 10kg/cm2."
$ws3.Range("A25").Value = "Providing and fixing ISI PVC Pipes includes jointing of pipes with one step PVC solvent cement, trenching, refilling & testing of joints complete as per direction of Engineer in Charge 32mm 10kg/cm2."
$ws2.Range("A22").Value = "Providing and fixing ISI PVC Pipes includes jointing of pipes with one step PVC solvent cement, trenching, refilling & testing of joints complete as per direction of Engineer in Charge 25mm 10kg/cm2."
$ws3.Range("A26").Value = "Providing and fixing ISI PVC Pipes includes jointing of pipes with one step PVC solvent cement, trenching, refilling & testing of joints complete as per direction of Engineer in Charge 25mm 10kg/cm2."
$ws2.Range("A23").Value = "Fabrication and Installation of Individual Hydrants with tap"
$ws3.Range("A27").Value = "Fabrication and Installation of Individual Hydrants with tap"
$ws2.Range("A24").Value = "Trenching for pipe laying"
$ws3.Range("A28").Value = "Trenching for pipe laying"
$ws2.Range("A25").Value = "Concrete Cutting for pipe laying in m3"
$ws3.Range("A29").Value = "Concrete Cutting for pipe laying in m3"
$ws2.Range("A26").Value = "Providing and laying in position cement concrete Filling 1:2:4 in m3"
$ws3.Range("A30").Value = "Providing and laying in position cement concrete Filling 1:2:4 in m3"
$ws3.Range("A3").Value = "Supply of Compressor Pump 1.5hp single phase ISI"
$ws3.Range("A4").Value = "Compressor Pump Fitting charge"
$ws3.Range("A5").Value = "GI pipe including providing, fitting, conveyance etc (With Trenching) 25mm pipe"
$ws3.Range("A6").Value = "Supply Fitting of HDPE ISI Pipe w/o data 20mm dia"
$ws3.Range("A7").Value = "Supply Fitting of ISI UPVC Pipe 32mm dia"
$ws3.Range("A8").Value = "Supply fitting of GI Union 32mm"
$ws3.Range("A9").Value = "Supply fitting of GI Bend 32mm"
$ws3.Range("A15").Value = "Fabrication and Installation of Metallic pump cover 1 x 0.8 x 0.9m"
$ws2.Range("A3").Value = "Supply Fitting of Submersible Pump Single Phase - 1.5hp 18 stage 107-30m head ISI"
$ws2.Range("A6").Value = "Supply Fitting of Submersible 3 Core Cable ISI 4mm"
$ws3.Range("A11").Value = "Supply Fitting of Submersible 3 Core Cable ISI 4mm"
$ws1.Range("A2").Value = "Earth work in excavation by mechanical means (Hydraulic excavator) /manual means in foundation trenches or drains (not exceeding 1.5 m in width or 10 sqm on plan), including dressing of sides and ramming of bottoms, lift up to 1.5 m, including getting out the excavated soil and disposal of surplus excavated soil as directed, within a lead of 50 m.`nAll kinds of soil (Code: 2.8.1)"
$ws1.Range("A3").Value = "Providing and laying in position cement concrete of specified grade excluding the cost of centering and shuttering - All work up to plinth level:`n1:2:4 (cement : 2 coarse sand : 4 graded stone aggregate 20mm nominal size) (Code: 4.1.3)"
$ws1.Range("A4").Value = "Providing and laying in position specified grade of reinforced cement concrete, excluding the cost of centering, shuttering, finishing and reinforcement - All work up to plinth level:`n1:2:4 ( 1 cement : 2 coarse sand : 4 graded stone aggregate 20mm nominal size) (Code: 5.1.3)"
$ws1.Range("A5").Value = "Steel reinforcement for R.C.C work including straightening, cutting, bending, placing in position and binding all complete upto plinth level`nMild steel and Medium Tensile steel bars (Code: 5.22.1)"
$ws1.Range("A6").Value = "Structural steel work in single section, fixed with or without connecting plate, including cutting, hoisting, fixing in position and applying a priming coat of approved steel primer all complete. (Code: 10.1)"
$ws1.Range("A7").Value = "Providing and fixing bolts including nuts and washers complete. (Code: 10.2)"

# Step 3: restore the numeric quantity values in column B.
$ws1.Range("B2").Value = 1
$ws1.Range("B3").Value = 2
$ws1.Range("B4").Value = 3
$ws1.Range("B5").Value = 4
$ws1.Range("B6").Value = 5
$ws1.Range("B7").Value = 6
$ws2.Range("B2").Value = 1
$ws2.Range("B3").Value = 1
$ws2.Range("B4").Value = 80
$ws2.Range("B5").Value = 85
$ws2.Range("B6").Value = 85
$ws2.Range("B7").Value = 1
$ws2.Range("B8").Value = 1
$ws2.Range("B9").Value = 1
$ws2.Range("B10").Value = 1
$ws2.Range("B11").Value = 1
$ws2.Range("B12").Value = 1
$ws2.Range("B13").Value = 1
$ws2.Range("B14").Value = 3000
$ws2.Range("B15").Value = 1
$ws2.Range("B16").Value = 1
$ws2.Range("B17").Value = 1
$ws2.Range("B18").Value = 1
$ws2.Range("B19").Value = 50
$ws2.Range("B20").Value = 50
$ws2.Range("B21").Value = 50
$ws2.Range("B22").Value = 50
$ws2.Range("B23").Value = 10
$ws2.Range("B24").Value = 100
$ws2.Range("B25").Value = 4
$ws2.Range("B26").Value = 4
$ws3.Range("B2").Value = 1
$ws3.Range("B3").Value = 1
$ws3.Range("B4").Value = 1
$ws3.Range("B5").Value = 15
$ws3.Range("B6").Value = 50
$ws3.Range("B7").Value = 60
$ws3.Range("B8").Value = 1
$ws3.Range("B9").Value = 1
$ws3.Range("B10").Value = 65
$ws3.Range("B11").Value = 20
$ws3.Range("B12").Value = 1
$ws3.Range("B13").Value = 1
$ws3.Range("B14").Value = 1
$ws3.Range("B15").Value = 1
$ws3.Range("B16").Value = 1
$ws3.Range("B17").Value = 1
$ws3.Range("B18").Value = 1
$ws3.Range("B19").Value = 3000
$ws3.Range("B20").Value = 1
$ws3.Range("B21").Value = 1
$ws3.Range("B22").Value = 1
$ws3.Range("B23").Value = 100
$ws3.Range("B24").Value = 100
$ws3.Range("B25").Value = 100
$ws3.Range("B26").Value = 100
$ws3.Range("B27").Value = 10
$ws3.Range("B28").Value = 250
$ws3.Range("B29").Value = 6
$ws3.Range("B30").Value = 6

# Step 4: restore sheet view state for the "Iron Structure" sheet (active sheet) -
# A2 becomes the selected cell.
$ws1.Range("A2").Select()
